$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Test Movie" (column G) is being updated so the UI no longer shows it as
# currently showing: its status changes from "Now Showing" to "End of Showing"
# and its age rating reverts from "PG13" to "PG".
$ws.Range("G4").Value = "End of Showing"
$ws.Range("G5").Value = "PG"

# Move the active selection to G4 to match the saved view state.
$ws.Range("G4").Select()
